$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new columns (student schema fields) that the note/progress
# report aggregation needs: technicalWritingApproved (before the old
# "backgroundPrepWorksheetApproved" column) and programProductRequirement
# (before the old "committeeCompApproved" column).
$ws.Columns("W").Insert()
$ws.Columns("AA").Insert()

# New column widths - copy the width from the column immediately to the
# left, matching Excel's default insert-column formatting behavior.
$ws.Columns("W").ColumnWidth = $ws.Columns("V").ColumnWidth
$ws.Columns("AA").ColumnWidth = $ws.Columns("Z").ColumnWidth

# Header row for the new fields
$ws.Range("W1").Value = "technicalWritingApproved"
$ws.Range("AA1").Value = "programProductRequirement"

# Sample row data for the new fields (dates, same format as neighboring
# approval-date columns)
$sampleDate = $ws.Range("T2").Value2
$ws.Range("W2").Value = $sampleDate
$ws.Range("AA2").Value = $sampleDate

# View state: scrolled over so the new columns are visible, with a
# different active selection
$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("Z11").Select()
